# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Gungnir_Profits workbook (per-class Leve profit sheets).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 7814325
$ws.Range("H40").Value = 5180763
$ws.Range("N40").Value = -1669697
$ws.Range("J40").Value = 1669347
$ws.Range("K40").Value = 7814325
$ws.Range("L40").Value = 1669347
$ws.Range("M40").Value = -7814150
$ws.Range("J43").Value = 1434.3158
$ws.Range("I43").Value = 71430840
$ws.Range("M43").Value = -71430771
$ws.Range("N43").Value = -1572.3158
$ws.Range("K43").Value = 71430840
$ws.Range("L43").Value = 1434.3158
$ws.Range("H43").Value = 19232428
$ws.Range("N138").Value = -19805
$ws.Range("K138").Value = 2020.4571
$ws.Range("L138").Value = 9525
$ws.Range("J138").Value = 3175
$ws.Range("H138").Value = 1458.2745
$ws.Range("I138").Value = 673.4857
$ws.Range("M138").Value = 3119.5429

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L5").Value = 100
$ws.Range("I5").Value = 260
$ws.Range("N5").Value = -324
$ws.Range("J5").Value = 100
$ws.Range("M5").Value = -148
$ws.Range("H5").Value = 180
$ws.Range("K5").Value = 260
$ws.Range("K63").Value = 2470.9048
$ws.Range("J63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("I63").Value = 2470.9048
$ws.Range("M63").Value = -1784.9048
$ws.Range("H63").Value = 2470.9048
$ws.Range("L63").Value = 0
$ws.Range("M66").Value = -8922.523999999999
$ws.Range("N66").ClearContents()
$ws.Range("I66").Value = 2470.9048
$ws.Range("H66").Value = 2470.9048
$ws.Range("L66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 12354.524

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -330
$ws.Range("J4").Value = 100
$ws.Range("M4").Value = -145
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 260
$ws.Range("K4").Value = 260
$ws.Range("H82").Value = 8915.571
$ws.Range("I82").Value = 3187.6667
$ws.Range("M82").Value = -2804.6667
$ws.Range("K82").Value = 3187.6667
$ws.Range("K85").Value = 3187.6667
$ws.Range("M85").Value = -1861.6667
$ws.Range("H85").Value = 8915.571
$ws.Range("I85").Value = 3187.6667
$ws.Range("H86").Value = 1554107
$ws.Range("L86").Value = 2328169.5
$ws.Range("I86").Value = 5982
$ws.Range("N86").Value = -2330415.5
$ws.Range("K86").Value = 5982
$ws.Range("M86").Value = -4859
$ws.Range("J86").Value = 2328169.5
$ws.Range("K89").Value = 29910
$ws.Range("J89").Value = 2328169.5
$ws.Range("M89").Value = -24294
$ws.Range("N89").Value = -11652079.5
$ws.Range("I89").Value = 5982
$ws.Range("L89").Value = 11640847.5
$ws.Range("H89").Value = 1554107
$ws.Range("L99").Value = 3130.1428
$ws.Range("I99").Value = 1010
$ws.Range("H99").Value = 1883
$ws.Range("M99").Value = 488
$ws.Range("J99").Value = 3130.1428
$ws.Range("N99").Value = -6126.1428
$ws.Range("K99").Value = 1010
$ws.Range("H116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("L118").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H122").Value = 0
$ws.Range("N123").Value = -52021.332
$ws.Range("J123").Value = 42221.332
$ws.Range("L123").Value = 42221.332
$ws.Range("H123").Value = 42221.332
$ws.Range("N125").ClearContents()
$ws.Range("J125").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("L125").Value = 0

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N6").Value = -2226
$ws.Range("L6").Value = 2000
$ws.Range("J6").Value = 2000
$ws.Range("H6").Value = 334066.66
$ws.Range("I7").Value = 48.285713
$ws.Range("K7").Value = 48.285713
$ws.Range("J7").Value = 20370
$ws.Range("M7").Value = 64.714287
$ws.Range("L7").Value = 20370
$ws.Range("H7").Value = 8515.666999999999
$ws.Range("N7").Value = -20596
$ws.Range("H16").Value = 1790.05
$ws.Range("N16").Value = -3874
$ws.Range("L16").Value = 3300
$ws.Range("I16").Value = 1142.9286
$ws.Range("K16").Value = 1142.9286
$ws.Range("M16").Value = -855.9286
$ws.Range("J16").Value = 3300
$ws.Range("J18").Value = 41000
$ws.Range("N18").Value = -41460
$ws.Range("H18").Value = 41000
$ws.Range("L18").Value = 41000
$ws.Range("H104").Value = 37742
$ws.Range("I104").Value = 0
$ws.Range("N104").Value = -42984
$ws.Range("L104").Value = 37742
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("J104").Value = 37742
$ws.Range("M105").Value = -40353
$ws.Range("L105").Value = 3895
$ws.Range("I105").Value = 42100
$ws.Range("K105").Value = 42100
$ws.Range("H105").Value = 25120
$ws.Range("N105").Value = -7389
$ws.Range("J105").Value = 3895
$ws.Range("N106").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("M107").Value = 1430.6111
$ws.Range("K107").Value = 489.3889
$ws.Range("N107").Value = -5537.4
$ws.Range("L107").Value = 1697.4
$ws.Range("J107").Value = 1697.4
$ws.Range("H107").Value = 920.8214
$ws.Range("I107").Value = 489.3889
$ws.Range("J109").Value = 26181.3
$ws.Range("L109").Value = 26181.3
$ws.Range("H109").Value = 26181.3
$ws.Range("N109").Value = -28261.3
$ws.Range("N111").Value = -43830.5
$ws.Range("L111").Value = 35650.5
$ws.Range("J111").Value = 35650.5
$ws.Range("H111").Value = 35650.5
$ws.Range("L112").Value = 875000
$ws.Range("J112").Value = 875000
$ws.Range("H112").Value = 875000
$ws.Range("N112").Value = -877954
$ws.Range("M113").Value = 1027.0714
$ws.Range("L113").Value = 3300
$ws.Range("J113").Value = 3300
$ws.Range("H113").Value = 1790.05
$ws.Range("I113").Value = 1142.9286
$ws.Range("N113").Value = -7640
$ws.Range("K113").Value = 1142.9286
$ws.Range("J114").Value = 35000
$ws.Range("L114").Value = 35000
$ws.Range("N114").Value = -43678
$ws.Range("H114").Value = 35000
$ws.Range("H115").Value = 47363.184
$ws.Range("I115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("L115").Value = 47363.184
$ws.Range("J115").Value = 47363.184
$ws.Range("N115").Value = -49713.184
$ws.Range("K115").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("N118").Value = -53114
$ws.Range("L118").Value = 49800
$ws.Range("H118").Value = 49800
$ws.Range("J118").Value = 49800
$ws.Range("J120").Value = 48650
$ws.Range("N120").Value = -55908
$ws.Range("L120").Value = 48650
$ws.Range("H120").Value = 48650
$ws.Range("J121").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("L121").Value = 0

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 38470624
$ws.Range("M132").Value = -5462
$ws.Range("I132").Value = 888
$ws.Range("K132").Value = 7992

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M113").Value = 1053.6666
$ws.Range("L113").Value = 1565
$ws.Range("J113").Value = 1565
$ws.Range("H113").Value = 1320.2727
$ws.Range("I113").Value = 1116.3334
$ws.Range("N113").Value = -5905
$ws.Range("K113").Value = 1116.3334

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K9").Value = 0
$ws.Range("J9").Value = 4742.2
$ws.Range("H9").Value = 4742.2
$ws.Range("M9").ClearContents()
$ws.Range("I9").Value = 0
$ws.Range("L9").Value = 4742.2
$ws.Range("N9").Value = -5190.2
$ws.Range("I18").Value = 1400
$ws.Range("K18").Value = 1400
$ws.Range("H18").Value = 1400
$ws.Range("M18").Value = -1228
$ws.Range("L61").Value = 1606.4
$ws.Range("N61").Value = -2010.4
$ws.Range("J61").Value = 1606.4
$ws.Range("K61").Value = 1540.8462
$ws.Range("I61").Value = 1540.8462
$ws.Range("M61").Value = -1338.8462
$ws.Range("H61").Value = 1569.3478
$ws.Range("M113").Value = 629.1538
$ws.Range("L113").Value = 1606.4
$ws.Range("J113").Value = 1606.4
$ws.Range("H113").Value = 1569.3478
$ws.Range("I113").Value = 1540.8462
$ws.Range("N113").Value = -5946.4
$ws.Range("K113").Value = 1540.8462

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M100").Value = -1095.3636
$ws.Range("K100").Value = 1636.3636
$ws.Range("I100").Value = 818.1818
$ws.Range("H100").Value = 884.8823
$ws.Range("L100").Value = 2014.3334
$ws.Range("N100").Value = -3096.3334
$ws.Range("J100").Value = 1007.1667
